# Auto-generated edit script: update cryptocurrency price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.138.65"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "2.508.77"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'571.36"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'166.06"
$ws.Range("E6").Value = "  -1.99%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.516"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").Value = "2.507.03"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "'0.353"
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "2.976.95"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "68.942.43"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "2.524.12"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").Value = "'7.61"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "'348.38"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'70.19"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").Value = "2.646.71"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "0.0₃0887"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").Value = "'7.82"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'460.30"
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'0.117"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").Value = "'157.82"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("D38").Value = "'19.00"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  -13.48%  "
$ws.Range("E46").Value = "  -6.31%  "
$ws.Range("D47").Value = "'141.28"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Value = "'0.526"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "'3.46"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("E51").Value = "  -3.42%  "
